# Auto-generated script applying scheduled market-data refresh to Maduin_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, matching the upstream data pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4313.8335
$ws.Range("J88").Value = 4313.8335
$ws.Range("L88").Value = 4313.8335
$ws.Range("N88").Value = -5125.8335
$ws.Range("H91").Value = 4313.8335
$ws.Range("J91").Value = 4313.8335
$ws.Range("L91").Value = 4313.8335
$ws.Range("N91").Value = -7121.8335
$ws.Range("H100").Value = 2493.1304
$ws.Range("I100").Value = 2711.7
$ws.Range("J100").Value = 1036
$ws.Range("K100").Value = 2711.7
$ws.Range("L100").Value = 1036
$ws.Range("M100").Value = -2170.7
$ws.Range("N100").Value = -2118
$ws.Range("H113").Value = 10099.833
$ws.Range("I113").Value = 2650
$ws.Range("K113").Value = 2650
$ws.Range("M113").Value = 604
$ws.Range("H116").Value = 6074.722
$ws.Range("I116").Value = 5386.9375
$ws.Range("K116").Value = 5386.9375
$ws.Range("M116").Value = -1944.9375
$ws.Range("H125").Value = 50003270
$ws.Range("I125").Value = 62501590
$ws.Range("K125").Value = 562514310
$ws.Range("M125").Value = -562511850
$ws.Range("H132").Value = 4044.1667
$ws.Range("I132").Value = 2599.4443
$ws.Range("J132").Value = 5488.8887
$ws.Range("K132").Value = 7798.3329
$ws.Range("L132").Value = 16466.6661
$ws.Range("M132").Value = -5268.3329
$ws.Range("N132").Value = -21526.6661
$ws.Range("H137").Value = 2186.3125
$ws.Range("I137").Value = 950.125
$ws.Range("J137").Value = 3422.5
$ws.Range("K137").Value = 2850.375
$ws.Range("L137").Value = 10267.5
$ws.Range("M137").Value = -300.375
$ws.Range("N137").Value = -15367.5
$ws.Range("H141").Value = 300
$ws.Range("I141").Value = 300
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 900
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 4280
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 366
$ws.Range("I11").Value = 299
$ws.Range("J11").Value = 399.5
$ws.Range("K11").Value = 299
$ws.Range("L11").Value = 399.5
$ws.Range("M11").Value = -155
$ws.Range("N11").Value = -687.5
$ws.Range("H45").Value = 1787.2858
$ws.Range("I45").Value = 1802.4
$ws.Range("J45").Value = 1749.5
$ws.Range("K45").Value = 1802.4
$ws.Range("L45").Value = 1749.5
$ws.Range("M45").Value = -1425.4
$ws.Range("N45").Value = -2503.5
$ws.Range("H61").Value = 2825.1538
$ws.Range("I61").Value = 2558
$ws.Range("K61").Value = 2558
$ws.Range("M61").Value = -2346
$ws.Range("H97").Value = 329.84616
$ws.Range("I97").Value = 303.41666
$ws.Range("J97").Value = 647
$ws.Range("K97").Value = 303.41666
$ws.Range("L97").Value = 647
$ws.Range("M97").Value = 192.58334
$ws.Range("N97").Value = -1639
$ws.Range("H102").Value = 2163
$ws.Range("I102").Value = 1112
$ws.Range("K102").Value = 1112
$ws.Range("M102").Value = 510
$ws.Range("H132").Value = 2349.5715
$ws.Range("I132").Value = 1241.1666
$ws.Range("K132").Value = 3723.4998
$ws.Range("M132").Value = -1193.4998
$ws.Range("H136").Value = 2825.1538
$ws.Range("I136").Value = 2558
$ws.Range("K136").Value = 7674
$ws.Range("M136").Value = -5124

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4576.375
$ws.Range("J105").Value = 5485.6665
$ws.Range("L105").Value = 5485.6665
$ws.Range("N105").Value = -8979.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2570.125
$ws.Range("I31").Value = 2330.077
$ws.Range("J31").Value = 3610.3333
$ws.Range("K31").Value = 2330.077
$ws.Range("L31").Value = 3610.3333
$ws.Range("M31").Value = -2035.077
$ws.Range("N31").Value = -4200.3333
$ws.Range("H34").Value = 2570.125
$ws.Range("I34").Value = 2330.077
$ws.Range("J34").Value = 3610.3333
$ws.Range("K34").Value = 2330.077
$ws.Range("L34").Value = 3610.3333
$ws.Range("M34").Value = -2128.077
$ws.Range("N34").Value = -4014.3333
$ws.Range("H58").Value = 996.3333
$ws.Range("I58").Value = 996.3333
$ws.Range("K58").Value = 996.3333
$ws.Range("M58").Value = -793.3333
$ws.Range("H94").Value = 3265.5833
$ws.Range("I94").Value = 3110.5715
$ws.Range("K94").Value = 3110.5715
$ws.Range("M94").Value = -2659.5715
$ws.Range("H99").Value = 6186.2354
$ws.Range("I99").Value = 5704.077
$ws.Range("K99").Value = 5704.077
$ws.Range("M99").Value = -4206.077
$ws.Range("H105").Value = 3869.1
$ws.Range("I105").Value = 1031.8334
$ws.Range("K105").Value = 1031.8334
$ws.Range("M105").Value = 715.1666
$ws.Range("H110").Value = 38999
$ws.Range("J110").Value = 38999
$ws.Range("L110").Value = 38999
$ws.Range("N110").Value = -47179
$ws.Range("H126").Value = 6186.2354
$ws.Range("I126").Value = 5704.077
$ws.Range("K126").Value = 17112.231
$ws.Range("M126").Value = -14642.231
$ws.Range("H132").Value = 5249.636
$ws.Range("I132").Value = 2356.5715
$ws.Range("J132").Value = 10312.5
$ws.Range("K132").Value = 7069.7145
$ws.Range("L132").Value = 30937.5
$ws.Range("M132").Value = -4539.7145
$ws.Range("N132").Value = -35997.5
$ws.Range("H134").Value = 2836.4119
$ws.Range("I134").Value = 2548.0667
$ws.Range("K134").Value = 7644.2001
$ws.Range("M134").Value = -5109.2001
$ws.Range("H136").Value = 996.3333
$ws.Range("I136").Value = 996.3333
$ws.Range("K136").Value = 2988.9999
$ws.Range("M136").Value = -438.9998999999998
$ws.Range("H141").Value = 392720.3
$ws.Range("J141").Value = 392720.3
$ws.Range("L141").Value = 392720.3
$ws.Range("N141").Value = -403080.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2900.5
$ws.Range("J39").Value = 2900.5
$ws.Range("L39").Value = 8701.5
$ws.Range("N39").Value = -9289.5
$ws.Range("H59").Value = 900
$ws.Range("I59").Value = 750
$ws.Range("J59").Value = 1200
$ws.Range("K59").Value = 2250
$ws.Range("L59").Value = 3600
$ws.Range("M59").Value = -1710
$ws.Range("N59").Value = -4680
$ws.Range("H98").Value = 596.25
$ws.Range("J98").Value = 630.75
$ws.Range("L98").Value = 1892.25
$ws.Range("N98").Value = -4888.25
$ws.Range("H122").Value = 1343.2858
$ws.Range("J122").Value = 1259.7778
$ws.Range("L122").Value = 11338.0002
$ws.Range("N122").Value = -16238.0002
$ws.Range("H132").Value = 2976.7693
$ws.Range("I132").Value = 1556.7142
$ws.Range("K132").Value = 14010.4278
$ws.Range("M132").Value = -11480.4278
$ws.Range("H139").Value = 4265.3335
$ws.Range("J139").Value = 4998
$ws.Range("L139").Value = 14994
$ws.Range("N139").Value = -25274

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 904.25
$ws.Range("I107").Value = 712.5
$ws.Range("J107").Value = 1096
$ws.Range("K107").Value = 712.5
$ws.Range("L107").Value = 1096
$ws.Range("M107").Value = 1207.5
$ws.Range("N107").Value = -4936

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3892.0527
$ws.Range("J46").Value = 4345.75
$ws.Range("L46").Value = 4345.75
$ws.Range("N46").Value = -4721.75
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H93").Value = 899.8182
$ws.Range("I93").Value = 899.7778
$ws.Range("K93").Value = 899.7778
$ws.Range("M93").Value = 348.2222
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
$ws.Range("H132").Value = 7497.5835
$ws.Range("I132").Value = 6871.375
$ws.Range("J132").Value = 8750
$ws.Range("K132").Value = 20614.125
$ws.Range("L132").Value = 26250
$ws.Range("M132").Value = -18084.125
$ws.Range("N132").Value = -31310
$ws.Range("H136").Value = 5310.5
$ws.Range("I136").Value = 5082.3335
$ws.Range("K136").Value = 15247.0005
$ws.Range("M136").Value = -12697.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1059.2
$ws.Range("I81").Value = 1074
$ws.Range("K81").Value = 2148
$ws.Range("M81").Value = -1087
$ws.Range("H84").Value = 1059.2
$ws.Range("I84").Value = 1074
$ws.Range("K84").Value = 10740
$ws.Range("M84").Value = -5436
